$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

# Column C ("Förändrad") holds a date serial number (45186 = 2023-09-17).
# Update every populated row (2 through 498) to the new date serial 45188 (2023-09-19).
$ws.Range("C2:C498").Value = 45188
